# Update the menu item names in the "Item Selection" column (A) to
# uppercase, fix the "Wine" duplicate -> "PIZZA" for the pasta-adjacent
# row, and shorten "Top Shelf" -> "Top Shf" in the "Purchase & Prep"
# column (C) for the spirits row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "STEAK"
$ws.Range("A12").Value = "CHICKEN"
$ws.Range("A14").Value = "PASTA"
$ws.Range("A16").Value = "PIZZA"
$ws.Range("A18").Value = "WINE"
$ws.Range("A20").Value = "SPIRITS"

$ws.Range("C20").Value = "Top Shf"

# Move the active selection from H20 to A20:C20.
[void]$ws.Range("A20:C20").Select()
